$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New rows of worked hours (Sprint 1 - MER)
$ws.Range("A29").Value = "Bruno Díaz"
$ws.Range("B29").Value = 42858
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = "Creación de Anteproyecto"
$ws.Range("E29").Value = "Arreglos"

$ws.Range("A30").Value = "Bruno Díaz"
$ws.Range("B30").Value = 42860
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = "Sprint 1 - MER"
$ws.Range("E30").Value = "Agregando atributos y Entidades al MER"

$ws.Range("A31").Value = "Bruno Díaz"
$ws.Range("B31").Value = 42861
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = "Sprint 1 - MER"
$ws.Range("E31").Value = "Modificando el MER a partir de la lista de requerimientos"

# Match the date number format used by the other date cells in column B
# (reuse the existing style via copy/paste-special instead of NumberFormat,
# which would otherwise mint a brand-new duplicate style entry)
$ws.Range("B28").Copy()
$ws.Range("B29:B31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view: scroll back to top and move the selection past the new data
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C32").Select()
